# Apply updated results from server for sheets 2025, 2030, 2035

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 0.001746599999999987
$ws.Range("E2").Value = 0.3768025662791765
$ws.Range("G2").Value = 0.2494892361374987
$ws.Range("I2").Value = 0.3498594026637402
$ws.Range("L2").Value = 0.6154130306695934
$ws.Range("M2").Value = 0.08155166666666669
$ws.Range("N2").Value = 12.74649432135023
$ws.Range("O2").Value = 3.479419494504655

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 0.06421966386250128
$ws.Range("E2").Value = 0.3707797323306978
$ws.Range("I2").Value = 0.4839130140029263
$ws.Range("L2").Value = 0.3654923693304067
$ws.Range("M2").Value = 0.08039441666666663
$ws.Range("N2").Value = 9.33903126345416
$ws.Range("O2").Value = 3.811842201374865

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.08566576916083338
$ws.Range("B2").Value = 0.02980781357332229
$ws.Range("E2").Value = 0.1657437090872776
$ws.Range("I2").Value = 0.4559350780436531
$ws.Range("M2").Value = 0.0437026666666667
$ws.Range("N2").Value = 8.958541033809858
$ws.Range("O2").Value = 5.258423268621357

$wb.Save()
